$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (single decimal point) need to be
# forced to text format first, otherwise Excel auto-converts them to numbers
# and trailing zeros / exact text formatting would be lost.
$textCells = @("D5", "D8", "D10", "D11", "D16", "D18", "D25", "D28", "D36", "D39", "D40", "D41", "D42", "D44", "D45", "D48", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.933.17"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.634.87"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "211.50"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "23.40"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "0.0610"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").Value = "0.0882"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.867.74"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "1.641.46"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "65.23"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "27.945.58"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "229.68"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("E19").Value = "  +3.61%  "
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "155.99"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").Value = "15.55"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D34").Value = "1.401.75"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").Value = "0.560"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "0.852"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("D44").Value = "66.12"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").Value = "5.45"
$ws.Range("D46").Value = "1.775.75"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").Value = "88.53"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "7.63"
$ws.Range("E51").Value = "  +1.83%  "
